$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Corrections to existing rows (columns F and/or G revised upward) ---
$updates = @{
    721 = @{ F = 27924 }
    723 = @{ F = 22560 }
    726 = @{ F = 35815 }
    727 = @{ F = 25095 }
    728 = @{ F = 24660 }
    729 = @{ F = 23225; G = 2517 }
    730 = @{ F = 19468; G = 2327 }
    732 = @{ F = 11839; G = 1900 }
    733 = @{ F = 31632; G = 3714 }
    734 = @{ F = 23080; G = 2540 }
    735 = @{ F = 19249; G = 2266 }
    736 = @{ F = 19475; G = 2180 }
    737 = @{ F = 18435; G = 2301 }
    738 = @{ F = 6484;  G = 933  }
    739 = @{ F = 8619;  G = 1384 }
    740 = @{ F = 24511; G = 2707 }
    741 = @{ F = 18814; G = 1904 }
    742 = @{ F = 17197; G = 1670 }
    743 = @{ F = 17914; G = 1591 }
    744 = @{ F = 14451; G = 1580 }
    745 = @{ F = 6102;  G = 908  }
    746 = @{ F = 7861;  G = 1214 }
    747 = @{ F = 21925; G = 2307 }
}

foreach ($r in $updates.Keys) {
    $row = $updates[$r]
    foreach ($col in $row.Keys) {
        $ws.Cells.Item([int]$r, $col).Value = $row[$col]
    }
}

# --- Row 748 previously had no F/G values; now populated ---
$ws.Cells.Item(748, "F").Value = 16659
$ws.Cells.Item(748, "G").Value = 1491

# --- New rows 749-753 appended with full data (A..G) ---
$newRows = @(
    @{ Row = 749; A = 44643; B = 1666482; C = 18588; D = 9649;  E = 19209; F = 14782; G = 1667 }
    @{ Row = 750; A = 44644; B = 1674586; C = 15364; D = 8104;  E = 19229; F = 14432; G = 1333 }
    @{ Row = 751; A = 44645; B = 1681865; C = 14726; D = 7279;  E = 19251; F = 9894;  G = 1158 }
    @{ Row = 752; A = 44646; B = 1687425; C = 11475; D = 5560;  E = 19270; F = 3131;  G = 484  }
    @{ Row = 753; A = 44647; B = 1690203; C = 5690;  D = 2778;  E = 19292; F = 3616;  G = 656  }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, "A").Value = $nr.A
    $ws.Cells.Item($r, "A").NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, "B").Value = $nr.B
    $ws.Cells.Item($r, "C").Value = $nr.C
    $ws.Cells.Item($r, "D").Value = $nr.D
    $ws.Cells.Item($r, "E").Value = $nr.E
    $ws.Cells.Item($r, "F").Value = $nr.F
    $ws.Cells.Item($r, "G").Value = $nr.G
}
